# Remove the grey dashed "separator" paragraphs and their associated
# small spacer paragraphs (an otherwise-empty paragraph with only
# w:spacing w:before="40") throughout the document, while leaving every
# other paragraph (including all inline images) untouched.
#
# Pattern being removed, wherever it occurs:
#   <w:p><w:pPr><w:spacing w:before="40"/></w:pPr></w:p>   (optional, not always present)
#   <w:p>
#     <w:pPr><w:spacing w:before="120" w:after="120"/></w:pPr>
#     <w:r><w:rPr><w:color w:val="CCCCCC"/><w:sz w:val="16"/></w:rPr>
#       <w:t>────────────────────────────────────────────────────────────</w:t>
#     </w:r>
#   </w:p>

$d = $word.ActiveDocument

$dashChar = [char]0x2500
$count = $d.Paragraphs.Count

for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    $trimmed = $text.Trim()

    $isSeparator = $trimmed.Length -gt 0 -and $trimmed.IndexOf($dashChar) -ge 0
    $isSpacer = ($trimmed.Length -eq 0) -and ([math]::Round($p.Format.SpaceBefore) -eq 2)

    if ($isSeparator -or $isSpacer) {
        $p.Range.Delete()
    }
}
